# "Generate Report for Handoff"
#
# The handoff report workbook gets regenerated: status text flips from the
# previous handback state to "ready for handoff", timestamps advance,
# translation priority flips from "ht" to "mt", and a new validation error
# (stale handback file) is recorded for the first file row on both the
# zh-cn and de-de sheets. A couple of report columns are also resized.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh       = $wb.Worksheets.Item("zh-cn")
$wsDe       = $wb.Worksheets.Item("de-de")

# ---- Status: "Handed back: in sync with en-US" -> "Ready for handoff" ----
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"

$wsZh.Range("C2").Value = "Ready for handoff"
$wsZh.Range("C3").Value = "Ready for handoff"

$wsDe.Range("C2").Value = "Ready for handoff"
$wsDe.Range("C3").Value = "Ready for handoff"

# ---- Latest HO Xliff Generate Date / de-de Latest Handoff Datetime ----
# "2016-10-18 13:18:54" -> "2016-10-18 13:20:32"
$wsOverview.Range("G2").Value = "2016-10-18 13:20:32"
$wsOverview.Range("G3").Value = "2016-10-18 13:20:32"

$wsDe.Range("H2").Value = "2016-10-18 13:20:32"
$wsDe.Range("H3").Value = "2016-10-18 13:20:32"

# ---- Priority: "ht" -> "mt" ----
$wsZh.Range("E2").Value = "mt"
$wsZh.Range("E3").Value = "mt"

$wsDe.Range("E2").Value = "mt"
$wsDe.Range("E3").Value = "mt"

# ---- zh-cn Latest Handoff Datetime: "2016-10-18 13:18:40" -> "2016-10-18 13:20:20" ----
$wsZh.Range("H2").Value = "2016-10-18 13:20:20"
$wsZh.Range("H3").Value = "2016-10-18 13:20:20"

# ---- New "Error Detail" message for the first file row on each language sheet ----
$errMsg = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a4c7ad83f6c875b629f52bdc803ac402485999da/e2e/afce9da6-ceeb-4336-ae7e-bd16b725c72e.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d897bb0cd4f320124137cffa194d5321115030d1/e2e/afce9da6-ceeb-4336-ae7e-bd16b725c72e.md."

$wsZh.Range("P2").Value = $errMsg
$wsDe.Range("P2").Value = $errMsg

# ---- Column width adjustments ----
# Overview!E:F and the zh-cn/de-de "Status"(C) column shrink from ~30 chars to ~17 chars;
# the zh-cn/de-de "Error Detail"(P) column grows from ~14 chars to 40 chars so the new
# long error message is readable.
$wsOverview.Range("E1").ColumnWidth = 16.3333333333333
$wsOverview.Range("F1").ColumnWidth = 16.3333333333333

$wsZh.Range("C1").ColumnWidth = 16.3333333333333
$wsZh.Range("P1").ColumnWidth = 39.1666666666667

$wsDe.Range("C1").ColumnWidth = 16.3333333333333
$wsDe.Range("P1").ColumnWidth = 39.1666666666667
